# Auto-generated: applies 2022-08-11 violent crime data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 4328
$ws.Range("I3").Value = 4531
$ws.Range("H4").Value = 1670
$ws.Range("I4").Value = 1042
$ws.Range("I5").Value = 411
$ws.Range("I6").Value = 4935
$ws.Range("H7").Value = 25981
$ws.Range("I7").Value = 15247

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 212
$ws.Range("I4").Value = 34
$ws.Range("I6").Value = 190
$ws.Range("I7").Value = 595

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I3").Value = 40
$ws.Range("I7").Value = 147

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 110
$ws.Range("I3").Value = 100
$ws.Range("I7").Value = 343

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 488
$ws.Range("I8").Value = 924
$ws.Range("I10").Value = 104
$ws.Range("I15").Value = 174
$ws.Range("I18").Value = 107
$ws.Range("I19").Value = 425
$ws.Range("I20").Value = 369
$ws.Range("I22").Value = 41
$ws.Range("I23").Value = 145
$ws.Range("I24").Value = 41
$ws.Range("I26").Value = 23
$ws.Range("I27").Value = 140
$ws.Range("I29").Value = 973
$ws.Range("I31").Value = 147
$ws.Range("I32").Value = 22
$ws.Range("I33").Value = 705
$ws.Range("I35").Value = 19
$ws.Range("I36").Value = 210
$ws.Range("I39").Value = 13
$ws.Range("I42").Value = 518
$ws.Range("I43").Value = 125
$ws.Range("I47").Value = 104
$ws.Range("I48").Value = 209
$ws.Range("I51").Value = 165
$ws.Range("H52").Value = 525
$ws.Range("I53").Value = 158
$ws.Range("I54").Value = 342
$ws.Range("I60").Value = 75
$ws.Range("I63").Value = 58
$ws.Range("I64").Value = 133
$ws.Range("I65").Value = 343
$ws.Range("I67").Value = 595
$ws.Range("I68").Value = 50
$ws.Range("I73").Value = 132
$ws.Range("I76").Value = 227
$ws.Range("I77").Value = 87
$ws.Range("I79").Value = 416
$ws.Range("I83").Value = 313
$ws.Range("I85").Value = 682
$ws.Range("I86").Value = 88
$ws.Range("I88").Value = 142
$ws.Range("I90").Value = 186
$ws.Range("I91").Value = 180
$ws.Range("I93").Value = 92
$ws.Range("I95").Value = 249
$ws.Range("I97").Value = 116
$ws.Range("H101").Value = 25981
$ws.Range("I101").Value = 15247

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 107
$ws.Range("I3").Value = 123
$ws.Range("I7").Value = 313

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 87
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 249

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 263
$ws.Range("I6").Value = 217
$ws.Range("I7").Value = 705

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 78
$ws.Range("I6").Value = 168
$ws.Range("I7").Value = 342

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 281
$ws.Range("I3").Value = 337
$ws.Range("I6").Value = 269
$ws.Range("I7").Value = 973

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I6").Value = 119
$ws.Range("I7").Value = 425

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I6").Value = 121
$ws.Range("I7").Value = 209

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 48
$ws.Range("I6").Value = 98
$ws.Range("I7").Value = 227

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I6").Value = 171
$ws.Range("I7").Value = 682

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 44
$ws.Range("I3").Value = 29
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 179
$ws.Range("I6").Value = 143
$ws.Range("I7").Value = 518

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 104

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("I2").Value = 15
$ws.Range("I7").Value = 41

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 50
$ws.Range("I7").Value = 145

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I2").Value = 57
$ws.Range("I7").Value = 180

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 133
$ws.Range("I6").Value = 120
$ws.Range("I7").Value = 416

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 38
$ws.Range("I7").Value = 133

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 101
$ws.Range("I3").Value = 115
$ws.Range("I7").Value = 369

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 107

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I6").Value = 63
$ws.Range("I7").Value = 210

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I3").Value = 24
$ws.Range("I7").Value = 92

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("H4").Value = 29
$ws.Range("H7").Value = 525

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I3").Value = 34
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 104

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 174

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 23

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("I2").Value = 2
$ws.Range("I6").Value = 13

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 19

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I3").Value = 42
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 132

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I6").Value = 69
$ws.Range("I7").Value = 116

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 39
$ws.Range("I3").Value = 51
$ws.Range("I6").Value = 42
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("I2").Value = 8
$ws.Range("I7").Value = 22

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 290
$ws.Range("I3").Value = 262
$ws.Range("I4").Value = 53
$ws.Range("I6").Value = 293
$ws.Range("I7").Value = 924

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 39
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 43
$ws.Range("I7").Value = 88

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I3").Value = 40
$ws.Range("I7").Value = 186

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 66
$ws.Range("I7").Value = 165

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("I3").Value = 15
$ws.Range("I7").Value = 50

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 125

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I5").Value = 5
$ws.Range("I7").Value = 158

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 41

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I2").Value = 27
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 162
$ws.Range("I3").Value = 155
$ws.Range("I6").Value = 123
$ws.Range("I7").Value = 488
